$wb = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item("About")

# Delete old rows 12 and 13 (two of the four lines of the old SCoHIbP note)
$about.Rows("12:13").Delete() | Out-Null

# Replace the (now 2-line) note text and drop its special formatting
$about.Range("A10").Value = "This variable is used to convert estimated avoided premature mortalities"
$about.Range("A10").ClearFormats() | Out-Null
$about.Range("A11").Value = "(calculated using data from HOIpTP) into a dollar amount."
$about.Range("A11").ClearFormats() | Out-Null

# The row that used to be the blank A14 (style s=1) is now A12; clear it completely so it's truly blank
$about.Range("A12").Clear() | Out-Null

$about.Activate() | Out-Null
$about.Range("I18").Select() | Out-Null

Write-Output "about sheet updated"
